{"js": "// Replace the 100 arithmetic-expression cells (5 columns x 20 rows) in the\n// single table with their new values, in document/row-major order. The\n// title paragraph (\"2025-08-15 Friday\") is left untouched.\nconst newValues = [\n  \"11-8=\",\n  \"70-54=\",\n  \"82-62=\",\n  \"16+2=\",\n  \"52-34=\",\n  \"49+0=\",\n  \"37+34=\",\n  \"55+40=\",\n  \"15+2=\",\n  \"46+24=\",\n  \"53-45=\",\n  \"36-23=\",\n  \"69-27=\",\n  \"42-2=\",\n  \"63-27=\",\n  \"89-47=\",\n  \"39-17=\",\n  \"94-47=\",\n  \"55-1=\",\n  \"0+52=\",\n  \"44+34=\",\n  \"46+25=\",\n  \"36+15=\",\n  \"36+58=\",\n  \"43+44=\",\n  \"71-50=\",\n  \"66-41=\",\n  \"50+4=\",\n  \"63-34=\",\n  \"63-38=\",\n  \"48+3=\",\n  \"3+67=\",\n  \"24+60=\",\n  \"42+4=\",\n  \"92-75=\",\n  \"83+0=\",\n  \"30+66=\",\n  \"48-4=\",\n  \"58-43=\",\n  \"87-31=\",\n  \"92-29=\",\n  \"37+4=\",\n  \"25+21=\",\n  \"40+48=\",\n  \"78-70=\",\n  \"33-9=\",\n  \"44+26=\",\n  \"13+81=\",\n  \"65-25=\",\n  \"26+3=\",\n  \"14+45=\",\n  \"28+1=\",\n  \"42+51=\",\n  \"10+3=\",\n  \"6+40=\",\n  \"51+39=\",\n  \"59-37=\",\n  \"49+46=\",\n  \"85-30=\",\n  \"89-17=\",\n  \"87-65=\",\n  \"10+67=\",\n  \"56-28=\",\n  \"3+73=\",\n  \"53+4=\",\n  \"42+0=\",\n  \"46+3=\",\n  \"90-77=\",\n  \"78-14=\",\n  \"10+6=\",\n  \"34-6=\",\n  \"54+34=\",\n  \"84+2=\",\n  \"80+1=\",\n  \"14+33=\",\n  \"11+3=\",\n  \"65+24=\",\n  \"8+74=\",\n  \"41-38=\",\n  \"33-14=\",\n  \"62-35=\",\n  \"46-9=\",\n  \"91+0=\",\n  \"77-46=\",\n  \"46+43=\",\n  \"86-74=\",\n  \"75-53=\",\n  \"63-7=\",\n  \"63-38=\",\n  \"36-34=\",\n  \"43+17=\",\n  \"43-20=\",\n  \"61-60=\",\n  \"83+8=\",\n  \"51-32=\",\n  \"88+2=\",\n  \"58-29=\",\n  \"65+33=\",\n  \"90-79=\",\n  \"94-15=\"\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = newValues.length / rowCount;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-expression cells (5 columns x 20 rows) in the\n# single table with their new values, in document/row-major order. The\n# title paragraph (\"2025-08-15 Friday\") is left untouched.\n$newValues = @(\n  \"11-8=\",\n  \"70-54=\",\n  \"82-62=\",\n  \"16+2=\",\n  \"52-34=\",\n  \"49+0=\",\n  \"37+34=\",\n  \"55+40=\",\n  \"15+2=\",\n  \"46+24=\",\n  \"53-45=\",\n  \"36-23=\",\n  \"69-27=\",\n  \"42-2=\",\n  \"63-27=\",\n  \"89-47=\",\n  \"39-17=\",\n  \"94-47=\",\n  \"55-1=\",\n  \"0+52=\",\n  \"44+34=\",\n  \"46+25=\",\n  \"36+15=\",\n  \"36+58=\",\n  \"43+44=\",\n  \"71-50=\",\n  \"66-41=\",\n  \"50+4=\",\n  \"63-34=\",\n  \"63-38=\",\n  \"48+3=\",\n  \"3+67=\",\n  \"24+60=\",\n  \"42+4=\",\n  \"92-75=\",\n  \"83+0=\",\n  \"30+66=\",\n  \"48-4=\",\n  \"58-43=\",\n  \"87-31=\",\n  \"92-29=\",\n  \"37+4=\",\n  \"25+21=\",\n  \"40+48=\",\n  \"78-70=\",\n  \"33-9=\",\n  \"44+26=\",\n  \"13+81=\",\n  \"65-25=\",\n  \"26+3=\",\n  \"14+45=\",\n  \"28+1=\",\n  \"42+51=\",\n  \"10+3=\",\n  \"6+40=\",\n  \"51+39=\",\n  \"59-37=\",\n  \"49+46=\",\n  \"85-30=\",\n  \"89-17=\",\n  \"87-65=\",\n  \"10+67=\",\n  \"56-28=\",\n  \"3+73=\",\n  \"53+4=\",\n  \"42+0=\",\n  \"46+3=\",\n  \"90-77=\",\n  \"78-14=\",\n  \"10+6=\",\n  \"34-6=\",\n  \"54+34=\",\n  \"84+2=\",\n  \"80+1=\",\n  \"14+33=\",\n  \"11+3=\",\n  \"65+24=\",\n  \"8+74=\",\n  \"41-38=\",\n  \"33-14=\",\n  \"62-35=\",\n  \"46-9=\",\n  \"91+0=\",\n  \"77-46=\",\n  \"46+43=\",\n  \"86-74=\",\n  \"75-53=\",\n  \"63-7=\",\n  \"63-38=\",\n  \"36-34=\",\n  \"43+17=\",\n  \"43-20=\",\n  \"61-60=\",\n  \"83+8=\",\n  \"51-32=\",\n  \"88+2=\",\n  \"58-29=\",\n  \"65+33=\",\n  \"90-79=\",\n  \"94-15=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
